$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# MaxTickersToProcess: 50 -> 20
$ws.Range("B9").Value = 20

# New config rows for batching the Yahoo Finance calls
$ws.Range("A13").Value = "BatchSize"
$ws.Range("B13").Value = 10

$ws.Range("A14").Value = "DelayBetweenBatchesMs"
$ws.Range("B14").Value = 700

# Match the author's final view/selection state
$ws.Range("B14").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
